$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Top servo housing and uprights adjusted" -- the servo-mount part costs
# (12-5V Converter line used for the top servo housing, row 8; Micro Servos
# line used for the uprights, row 9) were re-priced after the redesign.
# Unit Price column (C) drives the Total column (E) via a shared formula,
# so editing C8/C9 ripples through E8, E9 and the grand total in E11
# automatically on recalculation.
$ws.Range("C8").Value = 43
$ws.Range("C9").Value = 55

# Leave the cursor where the author left off after making the edit.
$ws.Range("D14").Select()
